# Actualización de flujos automatizados
# DatosRegistrarInformeVisitaVerificacion.xlsx
#
# Row 2 used to hold the verification-visit record for client "22296442".
# This record is replaced with a new client record ("21838047"). The cell
# is explicitly formatted as Text (so the numeric-looking id is preserved
# verbatim / keeps any leading zeros) before the value is written, and the
# data columns are (re)sized to fit their contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A2: new client code, stored as text -----------------------------------
# Setting the number format to Text ("@") *before* assigning the value makes
# Excel keep the numeric-looking string as text (shared string), matching
# the original authoring rather than re-interpreting it as a number.
$a2 = $ws.Range("A2")
$a2.NumberFormat = "@"
$a2.Value = "21838047"

# --- Column widths ----------------------------------------------------------
# Column B was widened manually; columns C:H were resized to fit their
# (unchanged) contents.
$ws.Columns("B").ColumnWidth = 13.833333333333332
$ws.Columns("C").ColumnWidth = 17.5
$ws.Columns("D").ColumnWidth = 12.0
$ws.Columns("E").ColumnWidth = 15.333333333333332
$ws.Columns("F").ColumnWidth = 13.166666666666668
$ws.Columns("G").ColumnWidth = 16.5
$ws.Columns("H").ColumnWidth = 9.0
